$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Φύλλο1")

# --- Table 1: Computation time (rows 2-4) ---
$ws.Range("B2").Value = 7.1673749999999998
$ws.Range("C2").Value = 7.1831829999999997
$ws.Range("D2").Value = 7.0804980000000004
$ws.Range("E2").Value = 7.2952789999999998

$ws.Range("B3").Value = 5.8552590000000002
$ws.Range("C3").Value = 5.7443160000000004
$ws.Range("D3").Value = 5.8854540000000002
$ws.Range("E3").Value = 5.7072919999999998

$ws.Range("B4").Value = 3.0548479999999998
$ws.Range("C4").Value = 3.0208490000000001
$ws.Range("D4").Value = 3.0833940000000002
$ws.Range("E4").Value = 2.996251

# --- Table 2: Overheads time (rows 10-12) ---
$ws.Range("B10").Value = 0.124359
$ws.Range("C10").Value = 0.127054
$ws.Range("D10").Value = 0.122822
$ws.Range("E10").Value = 0.122826

$ws.Range("B11").Value = 0.12756200000000001
$ws.Range("C11").Value = 0.14164499999999999
$ws.Range("D11").Value = 0.130415
$ws.Range("E11").Value = 0.14539299999999999

$ws.Range("B12").Value = 0.22983600000000001
$ws.Range("C12").Value = 0.22548199999999999
$ws.Range("D12").Value = 0.33199499999999998
$ws.Range("E12").Value = 0.238597

# --- Table 3: Total time (rows 18-20) ---
$ws.Range("B18").Value = 7.0430159999999997
$ws.Range("C18").Value = 7.0561290000000003
$ws.Range("D18").Value = 6.9576760000000002
$ws.Range("E18").Value = 7.1724540000000001

$ws.Range("B19").Value = 5.727697
$ws.Range("C19").Value = 5.602671
$ws.Range("D19").Value = 5.755039
$ws.Range("E19").Value = 5.5618990000000004

$ws.Range("B20").Value = 2.8250120000000001
$ws.Range("C20").Value = 2.7953670000000002
$ws.Range("D20").Value = 2.7513990000000002
$ws.Range("E20").Value = 2.7576529999999999

# Recalculate formulas (AVERAGE) and chart caches
$excel.CalculateFullRebuild() | Out-Null

# Update the selected cell in the sheet view
$ws.Range("E21").Select() | Out-Null
